$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

function Copy-Format($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------------------
# 1) Duplicate the whole original block (rows 1-10, A:S) sixteen rows down
#    (new block lives at rows 17-26), keeping formulas/relative refs intact.
# ---------------------------------------------------------------------------

# Row 1 (headers) -> Row 17
$ws.Range("A17").Value = $ws.Range("A1").Value()
$ws.Range("B17").Value = $ws.Range("B1").Value()
$ws.Range("C17").Value = $ws.Range("C1").Value()
$ws.Range("E17").Value = $ws.Range("E1").Value()
$ws.Range("F17").Value = $ws.Range("F1").Value()
$ws.Range("G17").Value = $ws.Range("G1").Value()
$ws.Range("I17").Value = $ws.Range("I1").Value()
$ws.Range("J17").Value = $ws.Range("J1").Value()
$ws.Range("K17").Value = $ws.Range("K1").Value()
$ws.Range("M17").Value = $ws.Range("M1").Value()
$ws.Range("N17").Value = $ws.Range("N1").Value()
$ws.Range("O17").Value = $ws.Range("O1").Value()
$ws.Range("Q17").Value = $ws.Range("Q1").Value()
$ws.Range("R17").Value = $ws.Range("R1").Value()
$ws.Range("S17").Value = $ws.Range("S1").Value()

# Row 2 (data row, group=0) -> Row 18
$ws.Range("A18").Value = 0
$ws.Range("B18").Value = 50491
$ws.Range("C18").Formula = "=B18/`$B`$20*100"
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 29079
$ws.Range("G18").Formula = "=F18/`$F`$20*100"
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 40830
$ws.Range("K18").Formula = "=J18/`$J`$20*100"
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = 22925
$ws.Range("O18").Formula = "=N18/`$N`$20*100"
$ws.Range("Q18").Value = 0
$ws.Range("R18").Value = 53363
$ws.Range("S18").Formula = "=R18/`$R`$20*100"

# Row 3 (data row, group=1) -> Row 19
$ws.Range("A19").Value = 1
$ws.Range("B19").Value = 4922
$ws.Range("C19").Formula = "=B19/`$B`$20*100"
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 17960
$ws.Range("G19").Formula = "=F19/`$F`$20*100"
$ws.Range("I19").Value = 1
$ws.Range("J19").Value = 13625
$ws.Range("K19").Formula = "=J19/`$J`$20*100"
$ws.Range("M19").Value = 1
$ws.Range("N19").Value = 36554
$ws.Range("O19").Formula = "=N19/`$N`$20*100"
$ws.Range("Q19").Value = 1
$ws.Range("R19").Value = 6846
$ws.Range("S19").Formula = "=R19/`$R`$20*100"

# Row 4 (sums/100%) -> Row 20
$ws.Range("B20").Formula = "=SUM(B18:B19)"
$ws.Range("C20").Formula = "=B20/`$B`$20*100"
$ws.Range("F20").Formula = "=SUM(F18:F19)"
$ws.Range("G20").Formula = "=F20/`$F`$20*100"
$ws.Range("J20").Formula = "=SUM(J18:J19)"
$ws.Range("K20").Formula = "=J20/`$J`$20*100"
$ws.Range("N20").Formula = "=SUM(N18:N19)"
$ws.Range("O20").Formula = "=N20/`$N`$20*100"
$ws.Range("R20").Formula = "=SUM(R18:R19)"
$ws.Range("S20").Formula = "=R20/`$R`$20*100"

# Row 5 (blank spacer w/ style) -> Row 21 (style copied below)

# Row 7 (total header) -> Row 23
$ws.Range("A23").Value = $ws.Range("A7").Value()
$ws.Range("B23").Value = $ws.Range("B7").Value()
$ws.Range("C23").Value = $ws.Range("C7").Value()

# Row 8 -> Row 24
$ws.Range("A24").Value = 1
$ws.Range("B24").Formula = "=B18+F18+J18+N18+R18"
$ws.Range("C24").Formula = "=B24/`$B`$26*100"

# Row 9 -> Row 25
$ws.Range("A25").Value = 0
$ws.Range("B25").Formula = "=B19+F19+J19+N19+R19"
$ws.Range("C25").Formula = "=B25/`$B`$26*100"

# Row 10 -> Row 26
$ws.Range("B26").Formula = "=B20+F20+J20+N20+R20"
$ws.Range("C26").Formula = "=B26/`$B`$26*100"

# ---------------------------------------------------------------------------
# Copy over the cell formatting (border/bold/alignment) from the originals
# onto the freshly duplicated block, cell by cell, so no stray blank/styled
# cells are introduced in unrelated columns.
# ---------------------------------------------------------------------------
Copy-Format "B1" "B17"
Copy-Format "F1" "F17"
Copy-Format "J1" "J17"
Copy-Format "N1" "N17"
Copy-Format "R1" "R17"

Copy-Format "A2" "A18"
Copy-Format "E2" "E18"
Copy-Format "I2" "I18"
Copy-Format "M2" "M18"
Copy-Format "Q2" "Q18"

Copy-Format "A3" "A19"
Copy-Format "E3" "E19"
Copy-Format "I3" "I19"
Copy-Format "M3" "M19"
Copy-Format "Q3" "Q19"

Copy-Format "B4" "B20"
Copy-Format "F4" "F20"
Copy-Format "J4" "J20"
Copy-Format "N4" "N20"
Copy-Format "R4" "R20"

Copy-Format "B5" "B21"

Copy-Format "B7" "B23"
Copy-Format "A8" "A24"
Copy-Format "A9" "A25"

# ---------------------------------------------------------------------------
# 2) Strip the "%" helper columns (C/G/K/O/S) out of the original block
#    (rows 1-4) - only Area/group columns remain there now.
# ---------------------------------------------------------------------------
$pctCols = @("C","G","K","O","S")
foreach ($col in $pctCols) {
    for ($r = 1; $r -le 4; $r++) {
        $ws.Range("$col$r").ClearContents()
    }
}

# ---------------------------------------------------------------------------
# 3) Insert a dark separator bar on row 15 (A:R) above the duplicated block.
# ---------------------------------------------------------------------------
$ws.Range("A15:R15").Interior.Color = 0

# ---------------------------------------------------------------------------
# 4) Restore the selection to match the saved view.
# ---------------------------------------------------------------------------
$ws.Range("F12").Select()
